$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(4, 9).Value = 'sv'
$ws.Cells.Item(4, 10).Value = 'Statement-opinion'
$ws.Cells.Item(8, 9).Value = 'sd'
$ws.Cells.Item(8, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(14, 9).Value = 'sv'
$ws.Cells.Item(14, 10).Value = 'Statement-opinion'
$ws.Cells.Item(17, 9).Value = 'sv'
$ws.Cells.Item(17, 10).Value = 'Statement-opinion'
$ws.Cells.Item(35, 9).Value = 'sv'
$ws.Cells.Item(35, 10).Value = 'Statement-opinion'
$ws.Cells.Item(48, 9).Value = 'sv'
$ws.Cells.Item(48, 10).Value = 'Statement-opinion'
$ws.Cells.Item(65, 9).Value = 'aa'
$ws.Cells.Item(65, 10).Value = 'Agree/Accept'
$ws.Cells.Item(66, 9).Value = 'sd'
$ws.Cells.Item(66, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(76, 9).Value = 'ba'
$ws.Cells.Item(76, 10).Value = 'Appreciation'
$ws.Cells.Item(81, 9).Value = 'sv'
$ws.Cells.Item(81, 10).Value = 'Statement-opinion'
$ws.Cells.Item(82, 9).Value = 'aa'
$ws.Cells.Item(82, 10).Value = 'Agree/Accept'
$ws.Cells.Item(139, 9).Value = 'aa'
$ws.Cells.Item(139, 10).Value = 'Agree/Accept'
$ws.Cells.Item(167, 9).Value = 'sv'
$ws.Cells.Item(167, 10).Value = 'Statement-opinion'
$ws.Cells.Item(173, 9).Value = 'sd'
$ws.Cells.Item(173, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(178, 9).Value = 'b'
$ws.Cells.Item(178, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(186, 9).Value = 'sd'
$ws.Cells.Item(186, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(205, 9).Value = 'sd'
$ws.Cells.Item(205, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(218, 9).Value = 'sd'
$ws.Cells.Item(218, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(221, 9).Value = 'sd'
$ws.Cells.Item(221, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(233, 9).Value = 'sd'
$ws.Cells.Item(233, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(236, 9).Value = 'sv'
$ws.Cells.Item(236, 10).Value = 'Statement-opinion'
$ws.Cells.Item(260, 9).Value = 'ba'
$ws.Cells.Item(260, 10).Value = 'Appreciation'
$ws.Cells.Item(289, 9).Value = 'sd'
$ws.Cells.Item(289, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(304, 9).Value = 'sd'
$ws.Cells.Item(304, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(329, 9).Value = 'sv'
$ws.Cells.Item(329, 10).Value = 'Statement-opinion'
$ws.Cells.Item(344, 9).Value = 'sd'
$ws.Cells.Item(344, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(347, 9).Value = 'b'
$ws.Cells.Item(347, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(348, 9).Value = 'sv'
$ws.Cells.Item(348, 10).Value = 'Statement-opinion'
$ws.Cells.Item(354, 9).Value = 'sv'
$ws.Cells.Item(354, 10).Value = 'Statement-opinion'
$ws.Cells.Item(355, 9).Value = 'sv'
$ws.Cells.Item(355, 10).Value = 'Statement-opinion'
